$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing PriceChange/UpDown values for row 31 (the last existing
# data row), which the repeater had left blank before this edit.
$ws.Range("X31").Value = -0.21000000000000085
$ws.Range("Y31").Value = "Down"

# Append the new data row (row 32) produced by running the repeater method
# against the 4 bio stocks.
$ws.Range("A32").Value = 42651.425115740742
$ws.Range("B32").Value = 13
$ws.Range("C32").Value = "Buy"
$ws.Range("D32").Value = 60
$ws.Range("E32").Value = 1075
$ws.Range("F32").Value = 181
$ws.Range("G32").Value = 71
$ws.Range("H32").Value = 28
$ws.Range("I32").Value = 95
$ws.Range("J32").Value = 4
$ws.Range("K32").Value = 13145
$ws.Range("L32").Value = 23
$ws.Range("M32").Value = 9
$ws.Range("N32").Value = 20
$ws.Range("O32").Value = 1
$ws.Range("P32").Value = "Noun"
$ws.Range("Q32").Value = 47.963765586266284
$ws.Range("R32").Value = 0.49
$ws.Range("S32").Value = 0.0521
$ws.Range("T32").Value = -0.0214
$ws.Range("U32").Value = 2.2799999999999998
$ws.Range("V32").Value = "N/A"
$ws.Range("W32").Value = 0

# Match the number formats used by the rest of the table: columns S/T are
# percentages (column A already inherits the column's date format).
$ws.Range("S32").NumberFormat = "0.00%"
$ws.Range("T32").NumberFormat = "0.00%"

# Refresh the "best fit" column widths now that new, wider data has landed.
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(2).ColumnWidth = 9.3333333333333333
$ws.Columns.Item(3).ColumnWidth = 9.6666666666666667
$ws.Columns.Item(4).ColumnWidth = 13.6666666666666667
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 13.6666666666666667
$ws.Columns.Item(7).ColumnWidth = 15.6666666666666667
$ws.Columns.Item(8).ColumnWidth = 15.6666666666666667
$ws.Columns.Item(9).ColumnWidth = 16.6666666666666667
$ws.Columns.Item(10).ColumnWidth = 16.6666666666666667
$ws.Columns.Item(12).ColumnWidth = 11.6666666666666667
$ws.Columns.Item(13).ColumnWidth = 11.6666666666666667
